$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 202; existing rows 202:247 shift down to 203:248,
# matching the growth of the table's dimension from R247 to R248.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row with the new weekly price record
# (same Mercado/Region/Categoria/etc. as surrounding rows, new date &
# price figures).
$ws.Cells.Item(202, 1).Value = 3
$ws.Cells.Item(202, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(202, 3).Value = "Coquimbo"
$ws.Cells.Item(202, 4).Value = 44943
$ws.Cells.Item(202, 5).Value = 5
$ws.Cells.Item(202, 6).Value = 100112010
$ws.Cells.Item(202, 7).Value = "Achicoria"
$ws.Cells.Item(202, 8).Value = "Sin especificar"
$ws.Cells.Item(202, 9).Value = "Primera"
$ws.Cells.Item(202, 10).Value = 100
$ws.Cells.Item(202, 11).Value = 7000
$ws.Cells.Item(202, 12).Value = 7500
$ws.Cells.Item(202, 13).Value = 7275
$ws.Cells.Item(202, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(202, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(202, 16).Value = 455
$ws.Cells.Item(202, 17).Value = 16
$ws.Cells.Item(202, 18).Value = "Hortaliza"
